$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: D-column "Price" values are stored as text (they use '.' as a
# thousands separator in some rows, e.g. "29.040.94"), so a leading apostrophe
# is used when assigning them to prevent Excel from reinterpreting them as
# numbers.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'29.040.94"
$ws.Range("E2").Value = "  -0.54%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.829.80"
$ws.Range("E3").Value = "  -0.22%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'241.17"
$ws.Range("E5").Value = "  -0.48%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.6231"
$ws.Range("E6").Value = "  -5.60%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "'0.07549"
$ws.Range("E8").Value = "  +1.80%  "

# Row 9 - OKB
$ws.Range("D9").Value = "'44.61"
$ws.Range("E9").Value = "  +6.60%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -0.65%  "

# Row 11 - Solana
$ws.Range("D11").Value = "'22.77"
$ws.Range("E11").Value = "  -0.47%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -1.89%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "'1.828.44"
$ws.Range("E13").Value = "  +0.17%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'4.959"
$ws.Range("E14").Value = "  -0.76%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.6647"
$ws.Range("E15").Value = "  -0.38%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'82.30"
$ws.Range("E16").Value = "  -0.69%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "'0.000009035"
$ws.Range("E17").Value = "  +7.55%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "'5.997"
$ws.Range("E18").Value = "  -1.85%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "'29.040.13"
$ws.Range("E19").Value = "  -0.47%  "

# Row 20 - WrappedliquidstakedEther2.0
$ws.Range("D20").Value = "'2.080.29"
$ws.Range("E20").Value = "  +0.59%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'224.78"
$ws.Range("E21").Value = "  -1.30%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "'12.34"
$ws.Range("E22").Value = "  -1.13%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "'7.209"
$ws.Range("E24").Value = "  +0.87%  "

# Row 25 - Dai
$ws.Range("D25").Value = "'1.000"
$ws.Range("E25").Value = "  +0.00%  "

# Row 26 - NEAR
$ws.Range("D26").Value = "'159.74"
$ws.Range("E26").Value = "  +0.50%  "

# Row 27
$ws.Range("D27").Value = "'8.388"
$ws.Range("E27").Value = "  -2.60%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -2.66%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "'17.84"
$ws.Range("E29").Value = "  -0.75%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "'1.495"
$ws.Range("E30").Value = "  -1.50%  "

# Row 31 - now Filecoin (was Toncoin)
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.051"
$ws.Range("E31").Value = "  -1.59%  "

# Row 32 - now Toncoin (was Filecoin)
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.205"
$ws.Range("E32").Value = "  +0.97%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'4.020"
$ws.Range("E33").Value = "  -0.75%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "'0.05215"
$ws.Range("E34").Value = "  -1.15%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "'1.839"
$ws.Range("E35").Value = "  -1.50%  "

# Row 36 - ARBITRUM
$ws.Range("E36").Value = "  +1.12%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "'0.7312"
$ws.Range("E37").Value = "  -1.67%  "

# Row 38 - HuobiToken
$ws.Range("E38").Value = "  -0.37%  "

# Row 39 - Maker
$ws.Range("D39").Value = "'1.276.73"
$ws.Range("E39").Value = "  -2.50%  "

# Row 40 - MXToken
$ws.Range("E40").Value = "  +0.51%  "

# Row 41 - VeChain
$ws.Range("D41").Value = "'0.01779"
$ws.Range("E41").Value = "  -0.81%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "'6.369"
$ws.Range("E42").Value = "  +7.13%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "'0.8905"
$ws.Range("E43").Value = "  -4.26%  "

# Row 44 - PaxDollar
$ws.Range("E44").Value = "  +0.08%  "

# Row 45 - Quant
$ws.Range("D45").Value = "'101.46"
$ws.Range("E45").Value = "  -1.30%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "'1.979.93"
$ws.Range("E46").Value = "  +0.46%  "

# Row 47 - Mantle
$ws.Range("D47").Value = "'0.5117"
$ws.Range("E47").Value = "  -0.54%  "

# Row 48 - Aave
$ws.Range("D48").Value = "'63.42"
$ws.Range("E48").Value = "  +0.65%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  -0.93%  "

# Row 50 - TheSandbox
$ws.Range("D50").Value = "'0.3962"
$ws.Range("E50").Value = "  -1.33%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "'8.878"
